# Edit: append a parenthetical "branch alternate" note (in red) to the
# first paragraph, and add a new shaded empty paragraph after the speech's
# closing paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document."
#    -> add two trailing spaces to the existing run, then append three
#       new red (C00000) runs spelling out the branch-alternate note.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$endNoMark = $p1.Range.End - 1   # exclude the paragraph mark

# Two trailing spaces stay part of the original (black) run.
$d.Range($endNoMark, $endNoMark).InsertAfter("  ")
$endNoMark = $endNoMark + 2

$redColor = 0x0000C0   # wdColor BGR for OOXML w:val="C00000"

$seg1 = "(This is a change " + [char]0x2013 + " Ve"
$d.Range($endNoMark, $endNoMark).InsertAfter($seg1)
$d.Range($endNoMark, $endNoMark + $seg1.Length).Font.Color = $redColor
$endNoMark = $endNoMark + $seg1.Length

$seg2 = "rsion for branch alternate"
$d.Range($endNoMark, $endNoMark).InsertAfter($seg2)
$d.Range($endNoMark, $endNoMark + $seg2.Length).Font.Color = $redColor
$endNoMark = $endNoMark + $seg2.Length

$seg3 = ")"
$d.Range($endNoMark, $endNoMark).InsertAfter($seg3)
$d.Range($endNoMark, $endNoMark + $seg3.Length).Font.Color = $redColor
$endNoMark = $endNoMark + $seg3.Length

# ---------------------------------------------------------------------
# 2) Add a new, empty, shaded paragraph right after the speech's final
#    paragraph ("... we are free at last.").
# ---------------------------------------------------------------------
$docEnd = $d.Content.End - 1
$null = $d.Range($docEnd, $docEnd).InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>')
